$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $val) {
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "64.402.36"
$ws.Range("E2").Value = "  -3.32%  "

Set-TextValue $ws.Range("D3") "3.137.24"
$ws.Range("E3").Value = "  -2.63%  "

$ws.Range("E4").Value = "  +0.24%  "

Set-TextValue $ws.Range("D5") "598.53"
$ws.Range("E5").Value = "  -1.70%  "

Set-TextValue $ws.Range("D6") "149.24"
$ws.Range("E6").Value = "  -5.95%  "

$ws.Range("E7").Value = "  +0.23%  "

Set-TextValue $ws.Range("D8") "3.130.34"
$ws.Range("E8").Value = "  -2.76%  "

Set-TextValue $ws.Range("D9") "0.533"
$ws.Range("E9").Value = "  -3.10%  "

Set-TextValue $ws.Range("D10") "0.153"
$ws.Range("E10").Value = "  -4.75%  "

Set-TextValue $ws.Range("D11") "5.61"
$ws.Range("E11").Value = "  -1.71%  "

Set-TextValue $ws.Range("D12") "0.477"
$ws.Range("E12").Value = "  -5.11%  "

Set-TextValue $ws.Range("D13") "0.0000259"
$ws.Range("E13").Value = "  -4.24%  "

Set-TextValue $ws.Range("D14") "36.90"
$ws.Range("E14").Value = "  -4.88%  "

Set-TextValue $ws.Range("D15") "3.662.07"
$ws.Range("E15").Value = "  -2.38%  "

Set-TextValue $ws.Range("D16") "64.600.06"
$ws.Range("E16").Value = "  -3.09%  "

$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextValue $ws.Range("D17") "3.161.61"
$ws.Range("E17").Value = "  -1.96%  "

$ws.Range("B18").Value = "TRON"
$ws.Range("C18").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
Set-TextValue $ws.Range("D18") "0.114"
$ws.Range("E18").Value = "  +0.35%  "

Set-TextValue $ws.Range("D19") "7.02"
$ws.Range("E19").Value = "  -4.59%  "

Set-TextValue $ws.Range("D20") "481.81"
$ws.Range("E20").Value = "  -4.96%  "

Set-TextValue $ws.Range("D21") "14.79"
$ws.Range("E21").Value = "  -2.32%  "

Set-TextValue $ws.Range("D22") "0.713"
$ws.Range("E22").Value = "  -2.70%  "

Set-TextValue $ws.Range("D23") "7.75"
$ws.Range("E23").Value = "  -3.07%  "

Set-TextValue $ws.Range("D24") "13.92"
$ws.Range("E24").Value = "  -4.64%  "

Set-TextValue $ws.Range("D25") "84.96"
$ws.Range("E25").Value = "  +0.13%  "

Set-TextValue $ws.Range("D26") "0.999"
$ws.Range("E26").Value = "  -0.18%  "

Set-TextValue $ws.Range("D27") "2.92"
$ws.Range("E27").Value = "  -2.71%  "

Set-TextValue $ws.Range("D28") "8.67"
$ws.Range("E28").Value = "  -4.97%  "

Set-TextValue $ws.Range("D29") "2.25"
$ws.Range("E29").Value = "  -4.55%  "

$ws.Range("B30").Value = "Hedera"
$ws.Range("C30").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue $ws.Range("D30") "0.122"
$ws.Range("E30").Value = "  +0.73%  "

$ws.Range("B31").Value = "NEARProtocol"
$ws.Range("C31").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextValue $ws.Range("D31") "7.12"
$ws.Range("E31").Value = "  +1.41%  "

$ws.Range("B32").Value = "FirstDigitalUSD"
$ws.Range("C32").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
Set-TextValue $ws.Range("D32") "1.00"
$ws.Range("E32").Value = "  -0.01%  "

$ws.Range("B33").Value = "Stacks"
$ws.Range("C33").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextValue $ws.Range("D33") "2.72"
$ws.Range("E33").Value = "  -7.24%  "

Set-TextValue $ws.Range("D34") "26.77"
$ws.Range("E34").Value = "  -4.79%  "

Set-TextValue $ws.Range("D35") "1.11"
$ws.Range("E35").Value = "  -6.58%  "

Set-TextValue $ws.Range("D36") "6.12"
$ws.Range("E36").Value = "  -5.31%  "

Set-TextValue $ws.Range("D37") "54.71"
$ws.Range("E37").Value = "  -1.40%  "

Set-TextValue $ws.Range("D38") "3.22"
$ws.Range("E38").Value = "  +4.35%  "

Set-TextValue $ws.Range("D39") "0.0₃0746"
$ws.Range("E39").Value = "  -3.46%  "

Set-TextValue $ws.Range("D40") "458.22"
$ws.Range("E40").Value = "  -8.52%  "

Set-TextValue $ws.Range("D41") "0.126"
$ws.Range("E41").Value = "  -4.62%  "

Set-TextValue $ws.Range("D42") "0.0401"
$ws.Range("E42").Value = "  -4.49%  "

Set-TextValue $ws.Range("D43") "8.52"
$ws.Range("E43").Value = "  -2.11%  "

Set-TextValue $ws.Range("D44") "2.42"
$ws.Range("E44").Value = "  -1.47%  "

Set-TextValue $ws.Range("D45") "2.887.26"
$ws.Range("E45").Value = "  -0.28%  "

Set-TextValue $ws.Range("D46") "0.274"
$ws.Range("E46").Value = "  -7.70%  "

Set-TextValue $ws.Range("D47") "26.91"
$ws.Range("E47").Value = "  -4.41%  "

$ws.Range("E48").Value = "  -0.05%  "

Set-TextValue $ws.Range("D49") "2.34"
$ws.Range("E49").Value = "  -3.23%  "

Set-TextValue $ws.Range("D50") "0.115"
$ws.Range("E50").Value = "  -0.24%  "

Set-TextValue $ws.Range("D51") "119.67"
$ws.Range("E51").Value = "  -2.07%  "
